$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from I1 into new J1, then set header texts
$ws.Cells.Item(1,9).Copy()
$ws.Cells.Item(1,10).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Headers: I1 becomes "Giornata 8", J1 becomes "Total average "
$ws.Cells.Item(1,9).Value = "Giornata 8"
$ws.Cells.Item(1,10).Value = "Total average "

# Data rows 2-21: I column becomes the new Giornata 8 scores,
# J column becomes the recomputed Total average over B:I (8 match days)
$ws.Cells.Item(2,9).Value = 5.785714285714286
$ws.Cells.Item(2,10).Value = 5.976219093406593
$ws.Cells.Item(3,9).Value = 6.615384615384615
$ws.Cells.Item(3,10).Value = 5.872856570512821
$ws.Cells.Item(4,9).Value = 6.576923076923077
$ws.Cells.Item(4,10).Value = 6.033012820512822
$ws.Cells.Item(5,9).Value = 5.5
$ws.Cells.Item(5,10).Value = 5.984775641025641
$ws.Cells.Item(6,9).Value = 5.692307692307693
$ws.Cells.Item(6,10).Value = 5.804601648351649
$ws.Cells.Item(7,9).Value = 5.333333333333333
$ws.Cells.Item(7,10).Value = 5.834612262737263
$ws.Cells.Item(8,9).Value = 5.642857142857143
$ws.Cells.Item(8,10).Value = 5.772744963369964
$ws.Cells.Item(9,9).Value = 6.066666666666666
$ws.Cells.Item(9,10).Value = 5.944432773109243
$ws.Cells.Item(10,9).Value = 6.576923076923077
$ws.Cells.Item(10,10).Value = 6.175730519480521
$ws.Cells.Item(11,9).Value = 6.428571428571429
$ws.Cells.Item(11,10).Value = 5.995386904761904
$ws.Cells.Item(12,9).Value = 6.730769230769231
$ws.Cells.Item(12,10).Value = 6.317055860805861
$ws.Cells.Item(13,9).Value = 5.730769230769231
$ws.Cells.Item(13,10).Value = 6.206267690642691
$ws.Cells.Item(14,9).Value = 5.071428571428571
$ws.Cells.Item(14,10).Value = 5.775869963369964
$ws.Cells.Item(15,9).Value = 6.833333333333333
$ws.Cells.Item(15,10).Value = 6.327953296703297
$ws.Cells.Item(16,9).Value = 5.576923076923077
$ws.Cells.Item(16,10).Value = 6.009354967948717
$ws.Cells.Item(17,9).Value = 6.0625
$ws.Cells.Item(17,10).Value = 6.110857371794872
$ws.Cells.Item(18,9).Value = 6.166666666666667
$ws.Cells.Item(18,10).Value = 5.908482142857142
$ws.Cells.Item(19,9).Value = 6
$ws.Cells.Item(19,10).Value = 5.793326465201465
$ws.Cells.Item(20,9).Value = 6.307692307692307
$ws.Cells.Item(20,10).Value = 5.962912087912088
$ws.Cells.Item(21,9).Value = 6.033333333333333
$ws.Cells.Item(21,10).Value = 6.120089285714285